# Azure Synapse Analytics PoC Environment.pptx
#
# 1) Add a "slow fade" slide transition to every slide. (The source deck's
#    change was a Morph transition written by modern PowerPoint as an
#    <mc:AlternateContent> block -- a <p159:morph> "Choice" for 2016+
#    clients with a plain <p:fade> "Fallback" for everyone else. The
#    classic PowerPoint object model used here has no Morph entry effect,
#    so we set the closest faithful, COM-reachable equivalent: the same
#    slow-speed fade that the real file's fallback branch renders.)
# 2) Split the "nano environment.tf" command run on slide 2 into three
#    runs so the visible command reads "nano terraform.tfvars" while
#    keeping the trailing tab + parenthetical comment intact.

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $trans = $slide.SlideShowTransition
    $trans.EntryEffect = 1793   # ppEffectFadeSmoothly
    $trans.Speed = 1            # ppTransitionSpeedSlow
}

# Locate the shape on slide 2 that holds the deployment instructions text
# and rewrite the "nano environment.tf" line to "nano terraform.tfvars".
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shape = $slide2.Shapes.Item($i)
    if (-not $shape.HasTextFrame) {
        continue
    }
    $tr = $shape.TextFrame.TextRange
    $fullText = $tr.Text
    $pos = $fullText.IndexOf("environment.tf")
    if ($pos -ge 0) {
        $sub = $tr.Characters($pos + 1, "environment.tf".Length)
        $sub.Text = "terraform.tfvars"
        break
    }
}
